# The document has a "title page" section layout, so each section exposes
# two distinct header/footer stories: the primary (Item 1) and the first-page
# (Item 2) header/footer. Both the Pearson logo (in the footers) and the
# BTEC logo (in the headers) appear once in each of those stories, and their
# embedded picture's OOXML "name" (wp:docPr/@name, mirrored onto
# pic:cNvPr/@name) needs to be swapped:
#   Pearson logo pictures: image1.png -> image2.png
#   BTEC logo pictures:    image2.jpg -> image1.jpg

$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {

    # Headers: BTec_Logo-Orange picture, rename image2.jpg -> image1.jpg
    # (the InlineShape.Name getter does not reflect the stored docPr/@name,
    # so match by the stable AlternativeText/descr instead of the old name)
    $hdrs = $sec.Headers
    for ($i = 1; $i -le $hdrs.Count; $i++) {
        $hdr = $hdrs.Item($i)
        if ($hdr.Exists) {
            $shapes = $hdr.Range.InlineShapes
            for ($j = 1; $j -le $shapes.Count; $j++) {
                $shp = $shapes.Item($j)
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp.Name = "image1.jpg"
                }
            }
        }
    }

    # Footers: PearsonLogo picture, rename image1.png -> image2.png
    $ftrs = $sec.Footers
    for ($i = 1; $i -le $ftrs.Count; $i++) {
        $ftr = $ftrs.Item($i)
        if ($ftr.Exists) {
            $shapes = $ftr.Range.InlineShapes
            for ($j = 1; $j -le $shapes.Count; $j++) {
                $shp = $shapes.Item($j)
                if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    $shp.Name = "image2.png"
                }
            }
        }
    }
}

Write-Output "Renamed header/footer logo pictures"
